$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.097.34"
$ws.Range("E2").Value = "  -2.47%  "
$ws.Range("D3").Value = "2.685.78"
$ws.Range("E3").Value = "  -2.84%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "549.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.66%  "
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("E8").Value = "  -2.34%  "
$ws.Range("E9").Value = "  -4.67%  "
$ws.Range("E10").Value = "  -2.97%  "
$ws.Range("E11").Value = "  -4.84%  "
$ws.Range("E12").Value = "  -11.93%  "
$ws.Range("D13").Value = "3.160.43"
$ws.Range("E13").Value = "  -2.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.15"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.53%  "
$ws.Range("D15").Value = "62.957.87"
$ws.Range("E15").Value = "  -1.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000147"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.83%  "
$ws.Range("D17").Value = "2.686.06"
$ws.Range("E17").Value = "  -3.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.95"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.93%  "
$ws.Range("E19").Value = "  -5.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "343.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.15%  "
$ws.Range("E21").Value = "  -5.56%  "
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.504"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.64"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.38%  "
$ws.Range("E25").Value = "  -1.79%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("E27").Value = "  -5.64%  "
$ws.Range("D28").Value = "0.0₃0859"
$ws.Range("E28").Value = "  -7.67%  "
$ws.Range("E29").Value = "  -2.49%  "
$ws.Range("E30").Value = "  -3.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "166.83"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.54%  "
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("E34").Value = "  -4.07%  "
$ws.Range("E35").Value = "  -3.45%  "
$ws.Range("E36").Value = "  -5.73%  "
$ws.Range("E37").Value = "  -3.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "339.10"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.91%  "
$ws.Range("E39").Value = "  -3.23%  "
$ws.Range("E40").Value = "  -7.37%  "
$ws.Range("E41").Value = "  -5.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.26"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.39"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.67%  "
$ws.Range("E44").Value = "  -8.12%  "
$ws.Range("E45").Value = "  -2.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0561"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.51%  "
$ws.Range("E47").Value = "  +0.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "11.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("E49").Value = "  -4.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "129.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.44%  "
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "2.093.58"
$ws.Range("E51").Value = "  -2.55%  "
